{"js": "// Fix initial template (ODM): the placeholder date paragraph reads\n// \"xx. x.  202\" (template remnants \"xx\"/\"x\" for day/month plus a\n// truncated placeholder year \"202\") and must become \"xx. x.  2025\" -\n// i.e. the missing \"5\" is appended so the placeholder year reads \"2025\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph holding the placeholder date - the one whose\n// text starts with the literal \"xx. x.\" template fragment.\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"xx. x.\") === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  // Search only within that paragraph for the truncated year \"202\" so we\n  // don't collide with unrelated \"20\" substrings elsewhere in the document\n  // (e.g. the student id \"F23209\").\n  const matches = target.search(\"202\", { matchCase: true, matchWholeWord: false });\n  matches.load(\"items\");\n  await context.sync();\n\n  if (matches.items.length > 0) {\n    const match = matches.items[0];\n    // Collapse to the end of the \"202\" match and insert the missing \"5\"\n    // right before it; inserting \"before\" the collapsed end point keeps\n    // the new character inside/merged with the matched run so it inherits\n    // the existing formatting (Lato-Regular font, yellow highlight, cs-CZ).\n    const endRange = match.getRange(\"End\");\n    endRange.insertText(\"5\", Word.InsertLocation.before);\n    await context.sync();\n  }\n}\n", "ps1": "# Fix initial template (ODM): the placeholder date paragraph \"xx. x.  202\"\n# (template remnants \"xx\"/\"x\" for day/month and a truncated year \"202\")\n# must become \"xx. x.  2025\" - i.e. the trailing year digits get the\n# missing \"5\" appended so the placeholder reads \"2025\".\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the placeholder date - it is the one\n# whose text starts with the literal \"xx. x.\" template fragment.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs($i)\n  if ($p.Range.Text.StartsWith(\"xx. x.\")) {\n    $target = $p\n    break\n  }\n}\n\nif ($target -ne $null) {\n  $r = $target.Range\n\n  # The paragraph range's End is exclusive of content and includes the\n  # hidden paragraph mark, so the last real character sits at End-2..End-1.\n  $lastCharRange = $d.Range($r.End - 2, $r.End - 1)\n\n  # Use Find/Replace scoped to just that trailing \"2\" so Word extends the\n  # existing (correctly formatted - Lato-Regular/yellow highlight/cs-CZ)\n  # run in place instead of fabricating a brand-new, unformatted run.\n  $find = $lastCharRange.Find\n  $find.ClearFormatting()\n  $find.Text = \"2\"\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = \"25\"\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nWrite-Output (\"final: [\" + $target.Range.Text + \"]\")\n"}
